$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J, matching the style used by the other
# header cells in row 1 (bold/centered/bordered style from B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J17
$data = @(
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(6, 7),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(4, 4),
    @(9, 9),
    @(4, 4),
    @(9, 9),
    @(9, 9),
    @(5, 6),
    @(9, 9),
    @(3, 3),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
